$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings: Score -> Intensity, Personal total -> Aspect intensity total
$ws.Range("J1").Value = "Sun Aspect Intensity"
$ws.Range("K1").Value = "Moon Aspect Intensity"
$ws.Range("L1").Value = "Asc Aspect Intensity"
$ws.Range("M1").Value = "Aspect intensity total"

# Add M-column rolling-sum formulas for rows 2..121 (each sums J:L of the row above)
for ($r = 2; $r -le 121; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 13).Formula = "=SUM(J$prev`:L$prev)"
}
